$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the "Family Medicine Preceptor Preference" (column H) entries for rows 3, 4, and 6
# these were incorrectly filled in with "Lakewood FHC" and should be blank based on
# location preferences.
$ws.Range("H3").ClearContents()
$ws.Range("H4").ClearContents()
$ws.Range("H6").ClearContents()

# Update the active selection to K4 to match the saved state after the edit
$ws.Range("K4").Select()
